# Orzeczenie_Tech.xlsx edit:
# - Add a new "Charakterystyka urzadzenia" (device characteristics) section
#   with 8 labelled fields (rows 11-15) to both the left (A:F) and right
#   (H:M) halves of the sheet.
# - Add a new "Opis Stanu Technicznego" (technical condition description)
#   section (rows 17-20) with a large merged free-text area, again
#   duplicated on both halves of the sheet.
# - Existing bordered cells (rows 7-9 template as well as the new rows)
#   get vertical="top" + wrapText="1" added to their alignment.
# - The "Zespol Orzekajacy:" block (previously rows 26/29/33) keeps its
#   texts, but the signature placeholder dot-lines get longer, and the
#   "Zatwierdzam" label shifts one column to the left (F->E, M->L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Merge the cells for the new rows first (so that the subsequent
#    format-paste gives every cell in the merge the exact same, uniform
#    bordered style - matching how the original template rows look).
# ---------------------------------------------------------------------

# Left block (columns A:F)
$ws.Range("A11:F11").Merge()
$ws.Range("A12:B12").Merge()
$ws.Range("A13:B13").Merge()
$ws.Range("A14:B14").Merge()
$ws.Range("A15:B15").Merge()
$ws.Range("D12:E12").Merge()
$ws.Range("D13:E13").Merge()
$ws.Range("D14:E14").Merge()
$ws.Range("D15:E15").Merge()
$ws.Range("A17:F17").Merge()
$ws.Range("A18:F20").Merge()

# Right block (columns H:M) - mirrors the left block
$ws.Range("H11:M11").Merge()
$ws.Range("H12:I12").Merge()
$ws.Range("H13:I13").Merge()
$ws.Range("H14:I14").Merge()
$ws.Range("H15:I15").Merge()
$ws.Range("K12:L12").Merge()
$ws.Range("K13:L13").Merge()
$ws.Range("K14:L14").Merge()
$ws.Range("K15:L15").Merge()
$ws.Range("H17:M17").Merge()
$ws.Range("H18:M20").Merge()

# ---------------------------------------------------------------------
# 2. Apply the same bordered formatting as the existing rows 7-9 to all
#    of the new cells (copy/paste-format keeps everything on one shared
#    style, exactly like the rest of the sheet).
# ---------------------------------------------------------------------

$ws.Range("A7:F7").Copy()
$ws.Range("A11:F11").PasteSpecial(-4122)
$ws.Range("A17:F20").PasteSpecial(-4122)

$ws.Range("A7:B7").Copy()
$ws.Range("A12:B15").PasteSpecial(-4122)

$ws.Range("C7:D7").Copy()
$ws.Range("D12:E15").PasteSpecial(-4122)

$ws.Range("H7:M7").Copy()
$ws.Range("H11:M11").PasteSpecial(-4122)
$ws.Range("H17:M20").PasteSpecial(-4122)

$ws.Range("H7:I7").Copy()
$ws.Range("H12:I15").PasteSpecial(-4122)

$ws.Range("J7:K7").Copy()
$ws.Range("K12:L15").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Fill in the text for the new "Charakterystyka urzadzenia" section.
# ---------------------------------------------------------------------

$ws.Range("A11").Value = "Charakterystyka urządzenia"
$ws.Range("H11").Value = "Charakterystyka urządzenia"

$ws.Range("A12").Value = "a) typ"
$ws.Range("H12").Value = "a) typ"
$ws.Range("D12").Value = "b) nr fab."
$ws.Range("K12").Value = "b) nr fab."

$ws.Range("A13").Value = "c) rok produkcji"
$ws.Range("H13").Value = "c) rok produkcji"
$ws.Range("D13").Value = "d) nr inw."
$ws.Range("K13").Value = "d) nr inw."

$ws.Range("A14").Value = "e) czas ekslpoatacji"
$ws.Range("H14").Value = "e) czas ekslpoatacji"
$ws.Range("D14").Value = "f) producent"
$ws.Range("K14").Value = "f) producent"

$ws.Range("A15").Value = "g) wartość księgowa"
$ws.Range("H15").Value = "g) wartość księgowa"
$ws.Range("D15").Value = "h) amortyzacja"
$ws.Range("K15").Value = "h) amortyzacja"

# ---------------------------------------------------------------------
# 4. Fill in the text for the new "Opis Stanu Technicznego" section.
# ---------------------------------------------------------------------

$ws.Range("A17").Value = "Opis Stanu Technicznego"
$ws.Range("H17").Value = "Opis Stanu Technicznego"

$placeholder = "skdjhngvuioabgfuavbrqberuovb[qeriogbf[aerioavbio[erh[aerigbre0hfwqgbkoer"
$ws.Range("A18").Value = $placeholder
$ws.Range("H18").Value = $placeholder

# ---------------------------------------------------------------------
# 5. Add vertical=top + wrap text to every bordered cell (template rows
#    7-9 as well as all the newly added rows 11-20) - this matches the
#    alignment change made to the shared cell style used by these
#    cells.
# ---------------------------------------------------------------------

$ws.Range("A7:F9").VerticalAlignment = -4160
$ws.Range("A7:F9").WrapText = $true
$ws.Range("H7:M9").VerticalAlignment = -4160
$ws.Range("H7:M9").WrapText = $true

$ws.Range("A11:F11").VerticalAlignment = -4160
$ws.Range("A11:F11").WrapText = $true
$ws.Range("A12:B15").VerticalAlignment = -4160
$ws.Range("A12:B15").WrapText = $true
$ws.Range("D12:E15").VerticalAlignment = -4160
$ws.Range("D12:E15").WrapText = $true
$ws.Range("A17:F20").VerticalAlignment = -4160
$ws.Range("A17:F20").WrapText = $true

$ws.Range("H11:M11").VerticalAlignment = -4160
$ws.Range("H11:M11").WrapText = $true
$ws.Range("H12:I15").VerticalAlignment = -4160
$ws.Range("H12:I15").WrapText = $true
$ws.Range("K12:L15").VerticalAlignment = -4160
$ws.Range("K12:L15").WrapText = $true
$ws.Range("H17:M20").VerticalAlignment = -4160
$ws.Range("H17:M20").WrapText = $true

# ---------------------------------------------------------------------
# 6. Update the "Zespol Orzekajacy" signature lines (longer dotted
#    lines) and move "Zatwierdzam" one column to the left.
# ---------------------------------------------------------------------

$ws.Range("A29").Value = "1 ........................."
$ws.Range("H29").Value = "1 ........................."
$ws.Range("A33").Value = "2 ........................."
$ws.Range("H33").Value = "2 ........................."

$ws.Range("F26").Value = ""
$ws.Range("M26").Value = ""
$ws.Range("E26").Value = "Zatwierdzam"
$ws.Range("L26").Value = "Zatwierdzam"
